# Correction of documentation for some entries of reason_for_exclusion:
# add a new "finished" column right before the existing last column
# (old column AE, "recoding_done", is pushed out to AF), and mark every
# data row as finished (value 1) in the new AE column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at AE (31st column), pushing the previous
#     AE ("recoding_done") one column to the right, into AF.
$ws.Columns.Item(31).Insert()
# The insert copies formatting from the column to the left (AD); strip
# it back out so the new column starts out unstyled, matching the rest
# of the freshly-added cells.
$ws.Columns.Item(31).ClearFormats()

# --- Header cell AE1 = "finished", using the same style as the other
#     bold/italic-ish header cells (W1..AB1 use style index 1).
$ws.Cells.Item(1, 23).Copy()
$ws.Cells.Item(1, 31).PasteSpecial(-4122)
$ws.Cells.Item(1, 31).Value = "finished"

# --- Every data row (2-101) gets a 1 in the new "finished" column.
for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 31).Value = 1
}

# --- Give the new column the same width as its neighbour (AD), which
#     previously was the only custom-width column in that area.
$ws.Columns.Item(31).ColumnWidth = $ws.Columns.Item(30).ColumnWidth

# --- Move the selection to reflect where editing ended up.
$ws.Range("AF101").Select() | Out-Null
